$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 410.33334
$ws.Range("I6").Value = 410.33334
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1231.00002
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1119.00002
$ws.Range("N6").Value = $null
$ws.Range("H19").Value = 2563.6667
$ws.Range("J19").Value = 1342
$ws.Range("L19").Value = 1342
$ws.Range("N19").Value = -1692
$ws.Range("H42").Value = 2041.4546
$ws.Range("I42").Value = 543.75
$ws.Range("K42").Value = 1631.25
$ws.Range("M42").Value = -1401.25
$ws.Range("H86").Value = 145399.58
$ws.Range("I86").Value = 202839.4
$ws.Range("J86").Value = 1800
$ws.Range("K86").Value = 202839.4
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -201716.4
$ws.Range("N86").Value = -4046
$ws.Range("H88").Value = 1812.375
$ws.Range("I88").Value = 2111
$ws.Range("K88").Value = 2111
$ws.Range("M88").Value = -1705
$ws.Range("H89").Value = 145399.58
$ws.Range("I89").Value = 202839.4
$ws.Range("J89").Value = 1800
$ws.Range("K89").Value = 1014197
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -1008581
$ws.Range("N89").Value = -20232
$ws.Range("H91").Value = 1812.375
$ws.Range("I91").Value = 2111
$ws.Range("K91").Value = 2111
$ws.Range("M91").Value = -707
$ws.Range("H92").Value = 126031.625
$ws.Range("J92").Value = 1994.25
$ws.Range("L92").Value = 1994.25
$ws.Range("N92").Value = -4490.25
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").Value = $null
$ws.Range("H107").Value = 1336.75
$ws.Range("I107").Value = 1204.5
$ws.Range("J107").Value = 1998
$ws.Range("K107").Value = 1204.5
$ws.Range("L107").Value = 1998
$ws.Range("M107").Value = 715.5
$ws.Range("N107").Value = -5838
$ws.Range("H115").Value = 1205.3636
$ws.Range("I115").Value = 433
$ws.Range("K115").Value = 1299
$ws.Range("M115").Value = 268
$ws.Range("H116").Value = 5477.154
$ws.Range("J116").Value = 5714.4
$ws.Range("L116").Value = 5714.4
$ws.Range("N116").Value = -12598.4
$ws.Range("H127").Value = 3099.5
$ws.Range("I127").Value = 3032.6667
$ws.Range("K127").Value = 9098.000100000001
$ws.Range("M127").Value = -4138.000100000001
$ws.Range("H129").Value = 934.7273
$ws.Range("I129").Value = 934.7273
$ws.Range("K129").Value = 2804.1819
$ws.Range("M129").Value = 2195.8181
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
$ws.Range("H131").Value = 1346.875
$ws.Range("I131").Value = 1346.875
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 4040.625
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 999.375
$ws.Range("N131").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4019.8572
$ws.Range("I61").Value = 4019.8572
$ws.Range("K61").Value = 4019.8572
$ws.Range("M61").Value = -3807.8572
$ws.Range("H74").Value = 51101.1
$ws.Range("I74").Value = 51101.1
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 51101.1
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -50227.1
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 51101.1
$ws.Range("I77").Value = 51101.1
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 255505.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -251137.5
$ws.Range("N77").Value = $null
$ws.Range("H107").Value = 149999.25
$ws.Range("J107").Value = 149999.25
$ws.Range("L107").Value = 149999.25
$ws.Range("N107").Value = -157679.25
$ws.Range("H110").Value = 2553.6667
$ws.Range("I110").Value = 2650.5
$ws.Range("K110").Value = 2650.5
$ws.Range("M110").Value = -605.5
$ws.Range("H132").Value = 17584.191
$ws.Range("I132").Value = 22280.584
$ws.Range("K132").Value = 66841.75199999999
$ws.Range("M132").Value = -64311.75199999999
$ws.Range("H136").Value = 4019.8572
$ws.Range("I136").Value = 4019.8572
$ws.Range("K136").Value = 12059.5716
$ws.Range("M136").Value = -9509.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1227
$ws.Range("I94").Value = 1116.6154
$ws.Range("K94").Value = 1116.6154
$ws.Range("M94").Value = -665.6153999999999
$ws.Range("H99").Value = 62416.94
$ws.Range("I99").Value = 86340.75
$ws.Range("J99").Value = 4999.8
$ws.Range("K99").Value = 86340.75
$ws.Range("L99").Value = 4999.8
$ws.Range("M99").Value = -84842.75
$ws.Range("N99").Value = -7995.8
$ws.Range("H107").Value = 2801
$ws.Range("I107").Value = 1603.6666
$ws.Range("K107").Value = 1603.6666
$ws.Range("M107").Value = 316.3334
$ws.Range("H112").Value = 149845
$ws.Range("J112").Value = 149845
$ws.Range("L112").Value = 149845
$ws.Range("N112").Value = -152799
$ws.Range("H134").Value = 2004.3636
$ws.Range("I134").Value = 1622.9524
$ws.Range("K134").Value = 4868.857199999999
$ws.Range("M134").Value = -2333.857199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2632
$ws.Range("I31").Value = 2119.9285
$ws.Range("J31").Value = 3348.9
$ws.Range("K31").Value = 2119.9285
$ws.Range("L31").Value = 3348.9
$ws.Range("M31").Value = -1824.9285
$ws.Range("N31").Value = -3938.9
$ws.Range("H34").Value = 2632
$ws.Range("I34").Value = 2119.9285
$ws.Range("J34").Value = 3348.9
$ws.Range("K34").Value = 2119.9285
$ws.Range("L34").Value = 3348.9
$ws.Range("M34").Value = -1917.9285
$ws.Range("N34").Value = -3752.9
$ws.Range("H104").Value = 97285
$ws.Range("J104").Value = 97285
$ws.Range("L104").Value = 97285
$ws.Range("N104").Value = -102527
$ws.Range("H122").Value = 1367.75
$ws.Range("I122").Value = 1296.3125
$ws.Range("K122").Value = 3888.9375
$ws.Range("M122").Value = -1438.9375
$ws.Range("H141").Value = 419982
$ws.Range("J141").Value = 419982
$ws.Range("L141").Value = 419982
$ws.Range("N141").Value = -430342

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1340697.4
$ws.Range("I4").Value = 925538.5
$ws.Range("J4").Value = 3001333
$ws.Range("K4").Value = 2776615.5
$ws.Range("L4").Value = 9003999
$ws.Range("M4").Value = -2776503.5
$ws.Range("N4").Value = -9004223
$ws.Range("H5").Value = 1327.625
$ws.Range("J5").Value = 1295.25
$ws.Range("L5").Value = 3885.75
$ws.Range("N5").Value = -4109.75
$ws.Range("H7").Value = 102.5
$ws.Range("I7").Value = 19
$ws.Range("J7").Value = 186
$ws.Range("K7").Value = 57
$ws.Range("L7").Value = 558
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = -782
$ws.Range("H68").Value = 600
$ws.Range("I68").Value = 600
$ws.Range("K68").Value = 1800
$ws.Range("M68").Value = -989
$ws.Range("H71").Value = 600
$ws.Range("I71").Value = 600
$ws.Range("K71").Value = 5400
$ws.Range("M71").Value = -1344
$ws.Range("H80").Value = 5331
$ws.Range("I80").Value = 3774
$ws.Range("J80").Value = 5850
$ws.Range("K80").Value = 11322
$ws.Range("L80").Value = 17550
$ws.Range("M80").Value = -10386
$ws.Range("N80").Value = -19422
$ws.Range("H83").Value = 5331
$ws.Range("I83").Value = 3774
$ws.Range("J83").Value = 5850
$ws.Range("K83").Value = 33966
$ws.Range("L83").Value = 52650
$ws.Range("M83").Value = -29286
$ws.Range("N83").Value = -62010
$ws.Range("H92").Value = 373.6087
$ws.Range("I92").Value = 359.93332
$ws.Range("J92").Value = 399.25
$ws.Range("K92").Value = 1079.79996
$ws.Range("L92").Value = 1197.75
$ws.Range("M92").Value = 168.2000400000002
$ws.Range("N92").Value = -3693.75
$ws.Range("H108").Value = 4700
$ws.Range("I108").Value = 4700
$ws.Range("K108").Value = 14100
$ws.Range("M108").Value = -11220
$ws.Range("H109").Value = 3474.5
$ws.Range("I109").Value = 3949
$ws.Range("J109").Value = 3000
$ws.Range("K109").Value = 11847
$ws.Range("L109").Value = 9000
$ws.Range("M109").Value = -10807
$ws.Range("N109").Value = -11080
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = $null
$ws.Range("N111").Value = $null
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = $null
$ws.Range("N112").Value = $null
$ws.Range("H113").Value = 809.7368
$ws.Range("J113").Value = 774.7273
$ws.Range("L113").Value = 2324.1819
$ws.Range("N113").Value = -6664.1819
$ws.Range("H114").Value = 753.4167
$ws.Range("I114").Value = 470.4
$ws.Range("J114").Value = 2168.5
$ws.Range("K114").Value = 1411.2
$ws.Range("L114").Value = 6505.5
$ws.Range("M114").Value = 1842.8
$ws.Range("N114").Value = -13013.5
$ws.Range("H116").Value = 795
$ws.Range("I116").Value = 795
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2385
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1057
$ws.Range("N116").Value = $null
$ws.Range("H117").Value = 5935
$ws.Range("I117").Value = 5935
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 17805
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = -14363
$ws.Range("N117").Value = $null
$ws.Range("H118").Value = 5988
$ws.Range("I118").Value = 5988
$ws.Range("K118").Value = 17964
$ws.Range("M118").Value = -16721
$ws.Range("H123").Value = 100030
$ws.Range("I123").Value = 100030
$ws.Range("K123").Value = 300090
$ws.Range("M123").Value = -297640
$ws.Range("H124").Value = 2115
$ws.Range("I124").Value = 930
$ws.Range("J124").Value = 3300
$ws.Range("K124").Value = 2790
$ws.Range("L124").Value = 9900
$ws.Range("M124").Value = 2120
$ws.Range("N124").Value = -19720
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").Value = $null
$ws.Range("H126").Value = 1330
$ws.Range("I126").Value = 1330
$ws.Range("K126").Value = 3990
$ws.Range("M126").Value = 950
$ws.Range("H128").Value = 148415.33
$ws.Range("I128").Value = 148415.33
$ws.Range("K128").Value = 445245.99
$ws.Range("M128").Value = -440265.99
$ws.Range("H129").Value = 11832.546
$ws.Range("I129").Value = 12515.8
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 37547.39999999999
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = -32547.39999999999
$ws.Range("N129").Value = -25000
$ws.Range("H130").Value = 9324.833000000001
$ws.Range("J130").Value = 10249
$ws.Range("L130").Value = 30747
$ws.Range("N130").Value = -40787
$ws.Range("H131").Value = 2134929.5
$ws.Range("I131").Value = 16478.857
$ws.Range("J131").Value = 2505658.2
$ws.Range("K131").Value = 49436.571
$ws.Range("L131").Value = 7516974.600000001
$ws.Range("M131").Value = -44396.571
$ws.Range("N131").Value = -7527054.600000001
$ws.Range("H132").Value = 1998.5
$ws.Range("I132").Value = 1998.5
$ws.Range("K132").Value = 17986.5
$ws.Range("M132").Value = -15456.5
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = $null
$ws.Range("N133").Value = $null
$ws.Range("H134").Value = 1800
$ws.Range("I134").Value = 1800
$ws.Range("K134").Value = 5400
$ws.Range("M134").Value = -330
$ws.Range("H135").Value = 1327.625
$ws.Range("J135").Value = 1295.25
$ws.Range("L135").Value = 11657.25
$ws.Range("N135").Value = -16727.25
$ws.Range("H137").Value = 5154.1
$ws.Range("J137").Value = 6000
$ws.Range("L137").Value = 18000
$ws.Range("N137").Value = -28200
$ws.Range("H138").Value = 1704
$ws.Range("I138").Value = 1482
$ws.Range("K138").Value = 4446
$ws.Range("M138").Value = 694
$ws.Range("H139").Value = 1818.3636
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null
$ws.Range("H140").Value = 3496.818
$ws.Range("I140").Value = 3083.5
$ws.Range("K140").Value = 9250.5
$ws.Range("M140").Value = -4070.5
$ws.Range("H141").Value = 1249.5
$ws.Range("I141").Value = 1249.5
$ws.Range("K141").Value = 3748.5
$ws.Range("M141").Value = 1431.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15098.8
$ws.Range("I70").Value = 14855.714
$ws.Range("K70").Value = 14855.714
$ws.Range("M70").Value = -14585.714
$ws.Range("H73").Value = 15098.8
$ws.Range("I73").Value = 14855.714
$ws.Range("K73").Value = 14855.714
$ws.Range("M73").Value = -13919.714
$ws.Range("H97").Value = 1056.4286
$ws.Range("I97").Value = 759.6667
$ws.Range("K97").Value = 759.6667
$ws.Range("M97").Value = -263.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = $null
$ws.Range("H46").Value = 7195.963
$ws.Range("I46").Value = 10530.9
$ws.Range("J46").Value = 5234.2354
$ws.Range("K46").Value = 10530.9
$ws.Range("L46").Value = 5234.2354
$ws.Range("M46").Value = -10342.9
$ws.Range("N46").Value = -5610.2354
$ws.Range("H70").Value = 53999.668
$ws.Range("J70").Value = 53999.668
$ws.Range("L70").Value = 53999.668
$ws.Range("N70").Value = -54539.668
$ws.Range("H73").Value = 53999.668
$ws.Range("J73").Value = 53999.668
$ws.Range("L73").Value = 53999.668
$ws.Range("N73").Value = -55871.668
$ws.Range("H122").Value = 2441.9492
$ws.Range("I122").Value = 1960
$ws.Range("J122").Value = 2591.889
$ws.Range("K122").Value = 5880
$ws.Range("L122").Value = 7775.667
$ws.Range("M122").Value = -3430
$ws.Range("N122").Value = -12675.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1625.0769
$ws.Range("I122").Value = 1629.6364
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 4888.9092
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -2438.9092
$ws.Range("N122").Value = -9700
$ws.Range("H132").Value = 60576.637
$ws.Range("I132").Value = 86305.64999999999
$ws.Range("J132").Value = 1399.9
$ws.Range("K132").Value = 258916.95
$ws.Range("L132").Value = 4199.700000000001
$ws.Range("M132").Value = -256386.95
$ws.Range("N132").Value = -9259.700000000001
